# Trade #42 closed at 2026-02-16 22:56:19 - base_strategy DOWN +0.000%
#
# A new trade log row (row 43) is appended to both the "All Trades" sheet
# and the per-strategy "base_strategy" sheet - the workbook mirrors every
# trade onto both tabs.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")
$row = 43

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item($row, 1).Value = 42                  # Trade #

    # Force the ISO date string to stay literal text instead of being
    # auto-converted to a date serial number by Excel's input parser.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"        # Date

    $ws.Cells.Item($row, 3).Value = "22:56:18"          # Time
    $ws.Cells.Item($row, 4).Value = "base_strategy"     # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"              # Side
    $ws.Cells.Item($row, 6).Value = 49.999998           # Entry Price
    $ws.Cells.Item($row, 7).Value = ""                  # Exit Price (still open)
    $ws.Cells.Item($row, 8).Value = "OPEN"              # Status
    $ws.Cells.Item($row, 9).Value = 0                   # P&L %
    $ws.Cells.Item($row, 10).Value = 0                  # P&L $
    $ws.Cells.Item($row, 11).Value = 100                # Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = ""                 # Exit Reason (still open)
    $ws.Cells.Item($row, 17).Value = 0                  # Duration (min)
}
